$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C width (value will be snapped by engine's internal rounding) ---
$ws.Columns.Item(3).ColumnWidth = 30.5

# --- Row 4: height change 45 -> 30 (content unchanged) ---
$ws.Rows.Item(4).RowHeight = 30

# --- Rows 15 & 16 grow to height 45 to fit the new wrapped text ---
$ws.Rows.Item(15).RowHeight = 45
$ws.Rows.Item(16).RowHeight = 45

# --- Row 8: replace text (adds a leading space before N) - new shared string #12 ---
$ws.Range("B8").Value = "simulate data with N(1, 4) for random intercept, MN for random slope and t3 for random error, called mn_t2_2re_30.Rdata"

# --- Row 15 C: new shared string #13 ---
$ws.Range("C15").Value = "two posterior data named are mn_t3_2re_10.Rdata and mn_t3_2re_new_10.Rdata"

# --- Row 11: replace text (adds a leading space before N) - new shared string #14 ---
$ws.Range("B11").Value = "simulate data with N(0, 4) for random intercept, MN for random slope and t3 for random error, called mn_t2_2re_30_new.Rdata"

# --- Row 12 C: new cell, new shared string #15 ---
$ws.Range("C12").Value = "use only the first 10 data sets"

# --- Row 15 B: new shared string #16 ---
$ws.Range("B15").Value = "Get results from previous simulation studies (using v3 for two data sets simulated above )and plot the densities to check the fit of the predictions"

# --- Row 17 C: new shared string #17 (plain) ---
$ws.Range("C17").Value = "create normal model file"

# --- Row 16 B: rich text shared string #18 ---
$ws.Range("B16").Value = "Run simulation on the complete 30 data sets using CDPM prior for BLQMM"
$ws.Range("B16").Characters(50, 5).Font.Bold = $true
$ws.Range("B16").Characters(55, 16).Font.Bold = $false

# --- Row 17 B: rich text shared string #19 ---
$ws.Range("B17").Value = "Run simulation on the complete 30 data sets using normal prior for BLQMM"
$ws.Range("B17").Characters(51, 12).Font.Bold = $true
$ws.Range("B17").Characters(63, 10).Font.Bold = $false

# --- Row 16 C: rich text shared string #20 ---
$ws.Range("C16").Value = "run additional 20 and combine with the previous results, run only on new data"
$ws.Range("C16").Characters(5, 13).Font.Underline = $true
$ws.Range("C16").Characters(18, 60).Font.Underline = $false

# --- Row 9: B/C swap values in place (adjust.. / us DPM..) already existing strings, order preserved ---
$ws.Range("B9").Value = "adjust the model file v2 to a new version v3 to incoperate the random slope"
$ws.Range("C9").Value = "us DPM to model both random effects"

# --- Row 11 C: x is also changed ... (existing string) ---
$ws.Range("C11").Value = "x is also changed to have mean 0"

# --- Row 12 B: use model_v3 ... (existing string) ---
$ws.Range("B12").Value = "use model_v3 to fit the data"

# --- Row 1 (header) values unchanged text, same shared strings, just reindexed automatically ---
$ws.Range("A1").Value = "DATE"
$ws.Range("B1").Value = "LOG"
$ws.Range("C1").Value = "NOTES"

# --- New date-log row 14 (copy format from the row 7 date-entry template) ---
$ws.Range("A7:C7").Copy($ws.Range("A14:C14"))
$ws.Range("A14").Value = 41884

# --- Selection change ---
$ws.Range("B20").Select()
